$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab) from "CopperA-HW45.xpc" to "CopperA"
$ws.Name = "CopperA"

# Add new row 16, mirroring the pattern (formatting) of row 15.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(16, 1).Value = 14

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
